$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (columns D, J, K, L, M, N, O, P, Q), derived from the
# weekly re-shuffle of the daily price records for Cilantro at the
# Macroferia Regional de Talca.
$rows = @{
    2  = @{ D = 44357; J = 150; K = 6500; L = 6500; M = 6500; N = '$/caja 20 docenas'; O = 'Región del Maule';      P = 6500; Q = 1  }
    3  = @{ D = 44362; J = 100; K = 6500; L = 6500; M = 6500; N = '$/caja 36 atados';  O = 'Región Metropolitana';  P = 181;  Q = 36 }
    4  = @{ D = 44371; J = 150; K = 6500; L = 6500; M = 6500; N = '$/caja 36 atados';  O = 'Región Metropolitana';  P = 181;  Q = 36 }
    5  = @{ D = 44355; J = 150; K = 7000; L = 7000; M = 7000; N = '$/caja 36 atados';  O = 'Región Metropolitana';  P = 194;  Q = 36 }
    6  = @{ D = 44369; J = 100; K = 7000; L = 7000; M = 7000; N = '$/caja 20 docenas'; O = 'Región Metropolitana';  P = 7000; Q = 1  }
    7  = @{ D = 44386; J = 200; K = 6500; L = 6500; M = 6500; N = '$/caja 36 atados';  O = 'Región Metropolitana';  P = 181;  Q = 36 }
    8  = @{ D = 44340; J = 150; K = 7000; L = 7000; M = 7000; N = '$/caja 36 atados';  O = 'Región del Maule';      P = 194;  Q = 36 }
    9  = @{ D = 44342; J = 150; K = 7000; L = 7000; M = 7000; N = '$/caja 36 atados';  O = 'Región del Maule';      P = 194;  Q = 36 }
    10 = @{ D = 44348; J = 150; K = 7000; L = 7000; M = 7000; N = '$/caja 36 atados';  O = 'Región del Maule';      P = 194;  Q = 36 }
    11 = @{ D = 44354; J = 150; K = 7000; L = 7000; M = 7000; N = '$/caja 36 atados';  O = 'Región del Maule';      P = 194;  Q = 36 }
    12 = @{ D = 44364; J = 100; K = 7000; L = 7000; M = 7000; N = '$/caja 36 atados';  O = 'Región Metropolitana';  P = 194;  Q = 36 }
    13 = @{ D = 44358; J = 150; K = 7000; L = 7000; M = 7000; N = '$/caja 36 atados';  O = 'Región Metropolitana';  P = 194;  Q = 36 }
    14 = @{ D = 44376; J = 150; K = 6500; L = 6500; M = 6500; N = '$/caja 36 atados';  O = 'Región Metropolitana';  P = 181;  Q = 36 }
    15 = @{ D = 44372; J = 150; K = 7000; L = 7000; M = 7000; N = '$/caja 36 atados';  O = 'Región Metropolitana';  P = 194;  Q = 36 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals.D   # D: Fecha
    $ws.Cells.Item($r, 10).Value = $vals.J   # J: Volumen
    $ws.Cells.Item($r, 11).Value = $vals.K   # K: Precio minimo
    $ws.Cells.Item($r, 12).Value = $vals.L   # L: Precio maximo
    $ws.Cells.Item($r, 13).Value = $vals.M   # M: Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $vals.N   # N: Unidad de comercializacion
    $ws.Cells.Item($r, 15).Value = $vals.O   # O: Origen
    $ws.Cells.Item($r, 16).Value = $vals.P   # P: Precio $/Kg
    $ws.Cells.Item($r, 17).Value = $vals.Q   # Q: Kg o Unidades
}
